$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The header in A1 ("stimuli") is renamed to "stimuli_0" to integrate the
# triggers into the experiment program. All other cell values are unchanged.
$ws.Range("A1").Value = "stimuli_0"

# Move the active selection, matching the saved cursor position in the
# workbook after the edit.
$ws.Range("D23").Select()
